$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.907.73"
$ws.Range("E2").Value = "  +4.99%  "

$ws.Range("D3").Value = "3.079.76"
$ws.Range("E3").Value = "  +4.19%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.12"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.57"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +8.51%  "

$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").Value = "3.072.26"
$ws.Range("E8").Value = "  +4.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.504"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +4.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.151"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.24"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.92%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +5.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000228"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +5.43%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.37"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +7.17%  "

$ws.Range("D15").Value = "3.579.99"
$ws.Range("E15").Value = "  +3.56%  "

$ws.Range("D16").Value = "63.931.43"
$ws.Range("E16").Value = "  +4.83%  "

$ws.Range("D17").Value = "3.079.62"
$ws.Range("E17").Value = "  +4.01%  "

$ws.Range("E18").Value = "  -0.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.78"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.62%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "486.58"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +7.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.61"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.689"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.23"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +7.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.99"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +5.91%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.62"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +8.41%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.77"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +5.97%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.05"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +6.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.01"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +11.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.20"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.16"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.01%  "

$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.44"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +9.67%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.82"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +9.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "55.70"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.02"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.65%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "473.95"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.95%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0824"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +7.21%  "

$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.176.43"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0400"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.66%  "

$ws.Range("E41").Value = "  +4.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.30"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "28.84"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +15.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.58"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +10.90%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.254"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +5.85%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.05"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +8.51%  "

$ws.Range("E48").Value = "  +3.94%  "

$ws.Range("D49").Value = "0.0₃0517"
$ws.Range("E49").Value = "  +3.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "116.85"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.59%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.08"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +7.86%  "
